# Updated results and code
# Applies the corrected stats values for rows 4-13 (stats_for_* metrics)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 0.296
$ws.Range("C4").Value = 0.051
$ws.Range("E4").Value = 0.151
$ws.Range("H4").Value = 0.189
$ws.Range("J4").Value = 0.108
$ws.Range("K4").Value = 0.342
$ws.Range("L4").Value = 0.102
$ws.Range("M4").Value = 0.32
$ws.Range("N4").Value = 0.271
$ws.Range("O4").Value = 0.02
$ws.Range("P4").Value = 0.143
$ws.Range("Q4").Value = 0.513
$ws.Range("R4").Value = 0.217
$ws.Range("S4").Value = 0.466
$ws.Range("T4").Value = 0.283
$ws.Range("W4").Value = 0.244
$ws.Range("Y4").Value = 0.208
$ws.Range("Z4").Value = 0.451
$ws.Range("AA4").Value = 0.133
$ws.Range("AB4").Value = 0.364
$ws.Range("AC4").Value = 0.126
$ws.Range("AE4").Value = 0.078
$ws.Range("AF4").Value = 0.737
$ws.Range("AG4").Value = 0.094
$ws.Range("AH4").Value = 0.307
$ws.Range("AI4").Value = 0.658
$ws.Range("AJ4").Value = 0.172
$ws.Range("AK4").Value = 0.415
$ws.Range("AL4").Value = 0.703
$ws.Range("AN4").Value = 0.339
$ws.Range("AO4").Value = 0.699

# Row 5
$ws.Range("B5").Value = 0.8159999999999999
$ws.Range("C5").Value = 0.15
$ws.Range("D5").Value = 0.388
$ws.Range("E5").Value = 0.658
$ws.Range("F5").Value = 0.225
$ws.Range("G5").Value = 0.474
$ws.Range("H5").Value = 0.8159999999999999
$ws.Range("I5").Value = 0.15
$ws.Range("J5").Value = 0.388
$ws.Range("K5").Value = 0.658
$ws.Range("L5").Value = 0.225
$ws.Range("M5").Value = 0.474
$ws.Range("N5").Value = 0.842
$ws.Range("O5").Value = 0.133
$ws.Range("P5").Value = 0.365
$ws.Range("Q5").Value = 0.579
$ws.Range("R5").Value = 0.244
$ws.Range("S5").Value = 0.494
$ws.Range("T5").Value = 0.579
$ws.Range("U5").Value = 0.244
$ws.Range("V5").Value = 0.494
$ws.Range("W5").Value = 0.737
$ws.Range("X5").Value = 0.194
$ws.Range("Y5").Value = 0.44
$ws.Range("Z5").Value = 0.8159999999999999
$ws.Range("AA5").Value = 0.15
$ws.Range("AB5").Value = 0.388
$ws.Range("AC5").Value = 0.763
$ws.Range("AD5").Value = 0.181
$ws.Range("AE5").Value = 0.425
$ws.Range("AF5").Value = 0.974
$ws.Range("AH5").Value = 0.16
$ws.Range("AI5").Value = 0.763
$ws.Range("AJ5").Value = 0.181
$ws.Range("AK5").Value = 0.425
$ws.Range("AL5").Value = 0.921
$ws.Range("AM5").Value = 0.073
$ws.Range("AN5").Value = 0.27
$ws.Range("AO5").Value = 0.886

# Row 6
$ws.Range("B6").Value = 0.434
$ws.Range("E6").Value = 0.246
$ws.Range("H6").Value = 0.307
$ws.Range("K6").Value = 0.45
$ws.Range("N6").Value = 0.41
$ws.Range("Q6").Value = 0.544
$ws.Range("T6").Value = 0.38
$ws.Range("W6").Value = 0.367
$ws.Range("Z6").Value = 0.581
$ws.Range("AC6").Value = 0.216
$ws.Range("AF6").Value = 0.839
$ws.Range("AI6").Value = 0.707
$ws.Range("AL6").Value = 0.797
$ws.Range("AO6").Value = 0.781

# Row 7
$ws.Range("B7").Value = 0.604
$ws.Range("E7").Value = 0.394
$ws.Range("H7").Value = 0.491
$ws.Range("K7").Value = 0.555
$ws.Range("N7").Value = 0.592
$ws.Range("Q7").Value = 0.5639999999999999
$ws.Range("T7").Value = 0.479
$ws.Range("W7").Value = 0.525
$ws.Range("Z7").Value = 0.702
$ws.Range("AC7").Value = 0.379
$ws.Range("AF7").Value = 0.915
$ws.Range("AI7").Value = 0.739
$ws.Range("AL7").Value = 0.867
$ws.Range("AO7").Value = 0.84

# Row 8
$ws.Range("B8").Value = 0.751
$ws.Range("C8").Value = 0.15
$ws.Range("D8").Value = 0.387
$ws.Range("E8").Value = 0.548
$ws.Range("F8").Value = 0.191
$ws.Range("G8").Value = 0.437
$ws.Range("H8").Value = 0.703
$ws.Range("I8").Value = 0.152
$ws.Range("J8").Value = 0.39
$ws.Range("K8").Value = 0.586
$ws.Range("L8").Value = 0.204
$ws.Range("M8").Value = 0.452
$ws.Range("N8").Value = 0.749
$ws.Range("O8").Value = 0.137
$ws.Range("P8").Value = 0.371
$ws.Range("Q8").Value = 0.55
$ws.Range("R8").Value = 0.229
$ws.Range("S8").Value = 0.479
$ws.Range("T8").Value = 0.501
$ws.Range("U8").Value = 0.206
$ws.Range("V8").Value = 0.454
$ws.Range("W8").Value = 0.662
$ws.Range("X8").Value = 0.18
$ws.Range("Y8").Value = 0.424
$ws.Range("Z8").Value = 0.751
$ws.Range("AA8").Value = 0.15
$ws.Range("AB8").Value = 0.387
$ws.Range("AC8").Value = 0.655
$ws.Range("AD8").Value = 0.177
$ws.Range("AE8").Value = 0.42
$ws.Range("AF8").Value = 0.893
$ws.Range("AG8").Value = 0.046
$ws.Range("AH8").Value = 0.215
$ws.Range("AI8").Value = 0.753
$ws.Range("AJ8").Value = 0.18
$ws.Range("AK8").Value = 0.424
$ws.Range("AL8").Value = 0.892
$ws.Range("AM8").Value = 0.078
$ws.Range("AN8").Value = 0.279
$ws.Range("AO8").Value = 0.846

# Row 9
$ws.Range("B9").Value = 0.658
$ws.Range("C9").Value = 0.225
$ws.Range("D9").Value = 0.474
$ws.Range("E9").Value = 0.421
$ws.Range("F9").Value = 0.244
$ws.Range("G9").Value = 0.494
$ws.Range("H9").Value = 0.579
$ws.Range("I9").Value = 0.244
$ws.Range("J9").Value = 0.494
$ws.Range("K9").Value = 0.5
$ws.Range("N9").Value = 0.632
$ws.Range("O9").Value = 0.233
$ws.Range("P9").Value = 0.482
$ws.Range("Q9").Value = 0.5
$ws.Range("T9").Value = 0.395
$ws.Range("U9").Value = 0.239
$ws.Range("V9").Value = 0.489
$ws.Range("W9").Value = 0.553
$ws.Range("X9").Value = 0.247
$ws.Range("Y9").Value = 0.497
$ws.Range("Z9").Value = 0.658
$ws.Range("AA9").Value = 0.225
$ws.Range("AB9").Value = 0.474
$ws.Range("AC9").Value = 0.553
$ws.Range("AD9").Value = 0.247
$ws.Range("AE9").Value = 0.497
$ws.Range("AF9").Value = 0.763
$ws.Range("AG9").Value = 0.181
$ws.Range("AH9").Value = 0.425
$ws.Range("AI9").Value = 0.737
$ws.Range("AJ9").Value = 0.194
$ws.Range("AK9").Value = 0.44
$ws.Range("AL9").Value = 0.842
$ws.Range("AM9").Value = 0.133
$ws.Range("AN9").Value = 0.365
$ws.Range("AO9").Value = 0.781

# Row 10
$ws.Range("B10").Value = 0.8159999999999999
$ws.Range("C10").Value = 0.15
$ws.Range("D10").Value = 0.388
$ws.Range("E10").Value = 0.579
$ws.Range("F10").Value = 0.244
$ws.Range("G10").Value = 0.494
$ws.Range("H10").Value = 0.737
$ws.Range("I10").Value = 0.194
$ws.Range("J10").Value = 0.44
$ws.Range("K10").Value = 0.658
$ws.Range("L10").Value = 0.225
$ws.Range("M10").Value = 0.474
$ws.Range("N10").Value = 0.8159999999999999
$ws.Range("O10").Value = 0.15
$ws.Range("P10").Value = 0.388
$ws.Range("Q10").Value = 0.579
$ws.Range("R10").Value = 0.244
$ws.Range("S10").Value = 0.494
$ws.Range("T10").Value = 0.579
$ws.Range("U10").Value = 0.244
$ws.Range("V10").Value = 0.494
$ws.Range("W10").Value = 0.737
$ws.Range("X10").Value = 0.194
$ws.Range("Y10").Value = 0.44
$ws.Range("Z10").Value = 0.8159999999999999
$ws.Range("AA10").Value = 0.15
$ws.Range("AB10").Value = 0.388
$ws.Range("AC10").Value = 0.658
$ws.Range("AD10").Value = 0.225
$ws.Range("AE10").Value = 0.474
$ws.Range("AF10").Value = 0.974
$ws.Range("AH10").Value = 0.16
$ws.Range("AI10").Value = 0.763
$ws.Range("AJ10").Value = 0.181
$ws.Range("AK10").Value = 0.425
$ws.Range("AL10").Value = 0.921
$ws.Range("AM10").Value = 0.073
$ws.Range("AN10").Value = 0.27
$ws.Range("AO10").Value = 0.886

# Row 11
$ws.Range("B11").Value = 0.8159999999999999
$ws.Range("C11").Value = 0.15
$ws.Range("D11").Value = 0.388
$ws.Range("E11").Value = 0.658
$ws.Range("F11").Value = 0.225
$ws.Range("G11").Value = 0.474
$ws.Range("H11").Value = 0.8159999999999999
$ws.Range("I11").Value = 0.15
$ws.Range("J11").Value = 0.388
$ws.Range("K11").Value = 0.658
$ws.Range("L11").Value = 0.225
$ws.Range("M11").Value = 0.474
$ws.Range("N11").Value = 0.842
$ws.Range("O11").Value = 0.133
$ws.Range("P11").Value = 0.365
$ws.Range("Q11").Value = 0.579
$ws.Range("R11").Value = 0.244
$ws.Range("S11").Value = 0.494
$ws.Range("T11").Value = 0.579
$ws.Range("U11").Value = 0.244
$ws.Range("V11").Value = 0.494
$ws.Range("W11").Value = 0.737
$ws.Range("X11").Value = 0.194
$ws.Range("Y11").Value = 0.44
$ws.Range("Z11").Value = 0.8159999999999999
$ws.Range("AA11").Value = 0.15
$ws.Range("AB11").Value = 0.388
$ws.Range("AC11").Value = 0.711
$ws.Range("AD11").Value = 0.206
$ws.Range("AE11").Value = 0.454
$ws.Range("AF11").Value = 0.974
$ws.Range("AH11").Value = 0.16
$ws.Range("AI11").Value = 0.763
$ws.Range("AJ11").Value = 0.181
$ws.Range("AK11").Value = 0.425
$ws.Range("AL11").Value = 0.921
$ws.Range("AM11").Value = 0.073
$ws.Range("AN11").Value = 0.27
$ws.Range("AO11").Value = 0.886

# Row 12
$ws.Range("B12").Value = 1.258
$ws.Range("C12").Value = 0.32
$ws.Range("D12").Value = 0.5659999999999999
$ws.Range("E12").Value = 1.68
$ws.Range("F12").Value = 1.098
$ws.Range("G12").Value = 1.048
$ws.Range("H12").Value = 1.613
$ws.Range("I12").Value = 1.334
$ws.Range("J12").Value = 1.155
$ws.Range("K12").Value = 1.4
$ws.Range("L12").Value = 0.5600000000000001
$ws.Range("M12").Value = 0.748
$ws.Range("N12").Value = 1.406
$ws.Range("O12").Value = 0.616
$ws.Range("P12").Value = 0.785
$ws.Range("Z12").Value = 1.258
$ws.Range("AA12").Value = 0.32
$ws.Range("AB12").Value = 0.5659999999999999
$ws.Range("AC12").Value = 1.793
$ws.Range("AD12").Value = 2.44
$ws.Range("AE12").Value = 1.562
$ws.Range("AF12").Value = 1.243
$ws.Range("AG12").Value = 0.238
$ws.Range("AH12").Value = 0.488
$ws.Range("AI12").Value = 1.034
$ws.Range("AJ12").Value = 0.033
$ws.Range("AK12").Value = 0.182
$ws.Range("AL12").Value = 1.086
$ws.Range("AM12").Value = 0.078
$ws.Range("AN12").Value = 0.28
$ws.Range("AO12").Value = 1.121

# Row 13
$ws.Range("B13").Value = 3.474
$ws.Range("C13").Value = 1.46
$ws.Range("D13").Value = 1.208
$ws.Range("E13").Value = 4.594
$ws.Range("F13").Value = 0.429
$ws.Range("G13").Value = 0.655
$ws.Range("H13").Value = 4.611
$ws.Range("I13").Value = 0.627
$ws.Range("J13").Value = 0.792
$ws.Range("K13").Value = 2.265
$ws.Range("L13").Value = 0.606
$ws.Range("M13").Value = 0.779
$ws.Range("N13").Value = 3.263
$ws.Range("O13").Value = 0.72
$ws.Range("P13").Value = 0.849
$ws.Range("Z13").Value = 2.514
$ws.Range("AA13").Value = 2.878
$ws.Range("AB13").Value = 1.697
$ws.Range("AC13").Value = 6.378
$ws.Range("AD13").Value = 2.181
$ws.Range("AE13").Value = 1.477
$ws.Range("AF13").Value = 1.605
$ws.Range("AG13").Value = 0.713
$ws.Range("AH13").Value = 0.844
$ws.Range("AI13").Value = 1.289
$ws.Range("AJ13").Value = 0.364
$ws.Range("AK13").Value = 0.603
$ws.Range("AL13").Value = 1.579
$ws.Range("AM13").Value = 0.717
$ws.Range("AN13").Value = 0.847
$ws.Range("AO13").Value = 1.491
